$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(10002, 3000121),
    @(10003, 3000122),
    @(10004, 3000123),
    @(10005, 3000124),
    @(10006, 3000125),
    @(10007, 3000126),
    @(10008, 3000127),
    @(10009, 3000128),
    @(10010, 3000129),
    @(10002, 3000130),
    @(10003, 3000131),
    @(10004, 3000132),
    @(10005, 3000133),
    @(10006, 3000134),
    @(10007, 3000135),
    @(10008, 3000136),
    @(10009, 3000137),
    @(10010, 3000138),
    @(10002, 3000139),
    @(10003, 3000140),
    @(10004, 3000141),
    @(10005, 3000142),
    @(10006, 3000143),
    @(10007, 3000144),
    @(10008, 3000145),
    @(10009, 3000146),
    @(10010, 3000147),
    @(10002, 3000148),
    @(10003, 3000149),
    @(10004, 3000150),
    @(10005, 3000151),
    @(10006, 3000152),
    @(10007, 3000153),
    @(10008, 3000154),
    @(10009, 3000155),
    @(10010, 3000156),
    @(10002, 3000157),
    @(10003, 3000158),
    @(10004, 3000159),
    @(10005, 3000160),
    @(10006, 3000161),
    @(10007, 3000162),
    @(10008, 3000163),
    @(10009, 3000164),
    @(10010, 3000165),
)

$startRow = 102
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $a = $data[$i][0]
    $b = $data[$i][1]
    $ws.Cells.Item($r, 1).Value = $a
    $ws.Cells.Item($r, 2).Value = $b
    $ws.Cells.Item($r, 3).Value = "eng"
    $ws.Cells.Item($r, 4).Value = $true
    $ws.Cells.Item($r, 5).Value = "superadmin"
    $ws.Cells.Item($r, 6).Value = "now()"
    $ws.Cells.Item($r, 7).Value = "now()"
}

# Match the author's final on-screen selection/scroll state after pasting
# the new rows (A102:B146 selected, A102 active).
$ws.Range("A102:B146").Select()
$excel.ActiveWindow.ScrollRow = 128

# Page setup: printed as portrait at 300 dpi.
$ws.PageSetup.Orientation = 1
$ws.PageSetup.HorizontalDpi = 300
$ws.PageSetup.VerticalDpi = 300
